$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill static (unchanging) columns for the 4 newly added rows (22-25)
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = 'Vega Modelo de Temuco'
$ws.Range("C22").Value = 'La Araucanía'
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 300000001
$ws.Range("G22").Value = 'Rabanito'
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("I22").Value = 'Primera'
$ws.Range("N22").Value = '$/docena de paquetes'
$ws.Range("Q22").Value = 12
$ws.Range("R22").Value = 'Hortaliza'

$ws.Range("A23").Value = 10
$ws.Range("B23").Value = 'Vega Modelo de Temuco'
$ws.Range("C23").Value = 'La Araucanía'
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 300000001
$ws.Range("G23").Value = 'Rabanito'
$ws.Range("H23").Value = 'Sin especificar'
$ws.Range("I23").Value = 'Primera'
$ws.Range("N23").Value = '$/docena de paquetes'
$ws.Range("Q23").Value = 12
$ws.Range("R23").Value = 'Hortaliza'

$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 'Vega Modelo de Temuco'
$ws.Range("C24").Value = 'La Araucanía'
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 300000001
$ws.Range("G24").Value = 'Rabanito'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("N24").Value = '$/docena de paquetes'
$ws.Range("Q24").Value = 12
$ws.Range("R24").Value = 'Hortaliza'

$ws.Range("A25").Value = 10
$ws.Range("B25").Value = 'Vega Modelo de Temuco'
$ws.Range("C25").Value = 'La Araucanía'
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 300000001
$ws.Range("G25").Value = 'Rabanito'
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("I25").Value = 'Primera'
$ws.Range("N25").Value = '$/docena de paquetes'
$ws.Range("Q25").Value = 12
$ws.Range("R25").Value = 'Hortaliza'

# Update / populate the per-row varying columns (D, J, K, L, M, O, P) for rows 7-25
$ws.Range("D7").Value = 44425
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("O7").Value = 'Provincia de Cautín'
$ws.Range("P7").Value = 583

$ws.Range("D8").Value = 44411
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("O8").Value = 'Provincia de Cautín'
$ws.Range("P8").Value = 583

$ws.Range("D9").Value = 44329
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5500
$ws.Range("O9").Value = 'Provincia de Cautín'
$ws.Range("P9").Value = 458

$ws.Range("D10").Value = 44424
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("O10").Value = 'Provincia de Cautín'
$ws.Range("P10").Value = 583

$ws.Range("D11").Value = 44166
$ws.Range("J11").Value = 55
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 6000
$ws.Range("O11").Value = 'Provincia de Cautín'
$ws.Range("P11").Value = 500

$ws.Range("D12").Value = 44299
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("O12").Value = 'Provincia de Cautín'
$ws.Range("P12").Value = 583

$ws.Range("D13").Value = 44162
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 5000
$ws.Range("O13").Value = 'Provincia de Cautín'
$ws.Range("P13").Value = 417

$ws.Range("D14").Value = 44427
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("O14").Value = 'Provincia de Cautín'
$ws.Range("P14").Value = 583

$ws.Range("D15").Value = 44413
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("O15").Value = 'Provincia de Cautín'
$ws.Range("P15").Value = 583

$ws.Range("D16").Value = 44410
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("O16").Value = 'Provincia de Cautín'
$ws.Range("P16").Value = 583

$ws.Range("D17").Value = 44327
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 6000
$ws.Range("O17").Value = 'Provincia de Cautín'
$ws.Range("P17").Value = 500

$ws.Range("D18").Value = 44196
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("O18").Value = 'Provincia de Cautín'
$ws.Range("P18").Value = 417

$ws.Range("D19").Value = 44369
$ws.Range("J19").Value = 20
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = 4000
$ws.Range("O19").Value = 'Región Metropolitana'
$ws.Range("P19").Value = 333

$ws.Range("D20").Value = 44195
$ws.Range("J20").Value = 55
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = 5000
$ws.Range("O20").Value = 'Provincia de Cautín'
$ws.Range("P20").Value = 417

$ws.Range("D21").Value = 44186
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 5000
$ws.Range("O21").Value = 'Provincia de Cautín'
$ws.Range("P21").Value = 417

$ws.Range("D22").Value = 44211
$ws.Range("J22").Value = 65
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 5000
$ws.Range("O22").Value = 'Provincia de Cautín'
$ws.Range("P22").Value = 417

$ws.Range("D23").Value = 44301
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("M23").Value = 6000
$ws.Range("O23").Value = 'Provincia de Cautín'
$ws.Range("P23").Value = 500

$ws.Range("D24").Value = 44326
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 6000
$ws.Range("M24").Value = 6000
$ws.Range("O24").Value = 'Provincia de Cautín'
$ws.Range("P24").Value = 500

$ws.Range("D25").Value = 44179
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = 6000
$ws.Range("O25").Value = 'Provincia de Cautín'
$ws.Range("P25").Value = 500

